# Apply the "exploratory" Shoe Size (D) and Eye Color (E) columns to the
# Data worksheet, and move the active selection the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Shoe Size (column D) and Eye Color (column E, L/G/R) for rows 2-15.
$shoeSize = @{
    2  = 8
    3  = 7
    4  = 6
    5  = 7.5
    6  = 9
    7  = 5.5
    8  = 9
    9  = 11
    10 = 5.5
    11 = 7
    12 = 10
    13 = 4.5
    14 = 5.5
    15 = 5
}

$eyeColor = @{
    2  = "L"
    3  = "G"
    4  = "G"
    5  = "R"
    6  = "G"
    7  = "G"
    8  = "R"
    9  = "R"
    10 = "R"
    11 = "L"
    12 = "L"
    13 = "R"
    14 = "G"
    15 = "L"
}

foreach ($row in 2..15) {
    $ws.Cells.Item($row, 4).Value = $shoeSize[$row]
    $ws.Cells.Item($row, 5).Value = $eyeColor[$row]
}

# Leave the selection where the author left it when they saved.
$ws.Range("E18").Select() | Out-Null
